$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values for rows 2-9 (columns G-T)
# Row 2
$ws.Range("G2").Value = 50.64491666666667
$ws.Range("H2").Value = 151.93475
$ws.Range("I2").Value = 0.3402395000245828
$ws.Range("J2").Value = 0.3402395000245828
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.5532856666666667
$ws.Range("N2").Value = 1.659857
$ws.Range("O2").Value = 0.4010144607159208
$ws.Range("P2").Value = 0.4010144607159208
$ws.Range("Q2").Value = 28.02110648119444
$ws.Range("R2").Value = 252.18995833075
$ws.Range("S2").Value = 0.1364409596166126
$ws.Range("T2").Value = 0.1364409596166126

# Row 3
$ws.Range("G3").Value = 50.64491666666667
$ws.Range("H3").Value = 151.93475
$ws.Range("I3").Value = 0.3402395000245828
$ws.Range("J3").Value = 0.3402395000245828
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.8264293333333333
$ws.Range("N3").Value = 2.479288
$ws.Range("O3").Value = 0.5989855392840792
$ws.Range("P3").Value = 0.5989855392840792
$ws.Range("Q3").Value = 41.85444471755556
$ws.Range("R3").Value = 376.690002458
$ws.Range("S3").Value = 0.2037985404079702
$ws.Range("T3").Value = 0.2037985404079702

# Row 4
$ws.Range("I4").Value = 0.5402707886290287
$ws.Range("J4").Value = 0.5402707886290287
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.5532856666666667
$ws.Range("N4").Value = 1.659857
$ws.Range("O4").Value = 0.4010144607159208
$ws.Range("P4").Value = 0.4010144607159208
$ws.Range("Q4").Value = 44.495084479489
$ws.Range("R4").Value = 400.455760315401
$ws.Range("S4").Value = 0.2166563989426352
$ws.Range("T4").Value = 0.2166563989426352

# Row 5
$ws.Range("I5").Value = 0.5402707886290287
$ws.Range("J5").Value = 0.5402707886290287
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.8264293333333333
$ws.Range("N5").Value = 2.479288
$ws.Range("O5").Value = 0.5989855392840792
$ws.Range("P5").Value = 0.5989855392840792
$ws.Range("Q5").Value = 66.46122467717601
$ws.Range("R5").Value = 598.1510220945839
$ws.Range("S5").Value = 0.3236143896863935
$ws.Range("T5").Value = 0.3236143896863935

# Row 6
$ws.Range("G6").Value = 17.514264
$ws.Range("H6").Value = 52.54279200000001
$ws.Range("I6").Value = 0.1176632289846506
$ws.Range("J6").Value = 0.1176632289846506
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.5532856666666667
$ws.Range("N6").Value = 1.659857
$ws.Range("O6").Value = 0.4010144607159208
$ws.Range("P6").Value = 0.4010144607159208
$ws.Range("Q6").Value = 9.690391233416001
$ws.Range("R6").Value = 87.21352110074402
$ws.Range("S6").Value = 0.04718465631737357
$ws.Range("T6").Value = 0.04718465631737358

# Row 7
$ws.Range("G7").Value = 17.514264
$ws.Range("H7").Value = 52.54279200000001
$ws.Range("I7").Value = 0.1176632289846506
$ws.Range("J7").Value = 0.1176632289846506
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.8264293333333333
$ws.Range("N7").Value = 2.479288
$ws.Range("O7").Value = 0.5989855392840792
$ws.Range("P7").Value = 0.5989855392840792
$ws.Range("Q7").Value = 14.474301521344
$ws.Range("R7").Value = 130.268713692096
$ws.Range("S7").Value = 0.07047857266727706
$ws.Range("T7").Value = 0.07047857266727706

# Row 8
$ws.Range("G8").Value = 0.2718733333333334
$ws.Range("H8").Value = 0.81562
$ws.Range("I8").Value = 0.001826482361737853
$ws.Range("J8").Value = 0.001826482361737852
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.5532856666666667
$ws.Range("N8").Value = 1.659857
$ws.Range("O8").Value = 0.4010144607159208
$ws.Range("P8").Value = 0.4010144607159208
$ws.Range("Q8").Value = 0.1504236184822222
$ws.Range("R8").Value = 1.35381256634
$ws.Range("S8").Value = 0.0007324458392994463
$ws.Range("T8").Value = 0.0007324458392994463

# Row 9
$ws.Range("G9").Value = 0.2718733333333334
$ws.Range("H9").Value = 0.81562
$ws.Range("I9").Value = 0.001826482361737853
$ws.Range("J9").Value = 0.001826482361737852
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.8264293333333333
$ws.Range("N9").Value = 2.479288
$ws.Range("O9").Value = 0.5989855392840792
$ws.Range("P9").Value = 0.5989855392840792
$ws.Range("Q9").Value = 0.2246840976177778
$ws.Range("R9").Value = 2.02215687856
$ws.Range("S9").Value = 0.001094036522438406
$ws.Range("T9").Value = 0.001094036522438406
